$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue 'D2' '59.353.88'
$ws.Range('E2').Value = '  -0.40%  '
Set-TextValue 'D3' '2.643.15'
$ws.Range('E3').Value = '  -0.07%  '
$ws.Range('E4').Value = '  -0.10%  '
Set-TextValue 'D5' '517.98'
$ws.Range('E5').Value = '  +0.00%  '
Set-TextValue 'D6' '146.23'
$ws.Range('E6').Value = '  -0.35%  '
$ws.Range('E7').Value = '  +0.23%  '
Set-TextValue 'D8' '0.574'
$ws.Range('E8').Value = '  +0.30%  '
Set-TextValue 'D9' '2.650.71'
$ws.Range('E9').Value = '  -0.72%  '
Set-TextValue 'D10' '6.31'
$ws.Range('E10').Value = '  -2.93%  '
$ws.Range('E11').Value = '  -1.32%  '
$ws.Range('E12').Value = '  -0.99%  '
$ws.Range('E13').Value = '  +0.75%  '
Set-TextValue 'D14' '3.106.45'
$ws.Range('E14').Value = '  +0.04%  '
Set-TextValue 'D15' '59.359.02'
$ws.Range('E15').Value = '  -0.19%  '
Set-TextValue 'D16' '21.05'
$ws.Range('E16').Value = '  -1.05%  '
Set-TextValue 'D18' '2.644.60'
$ws.Range('E18').Value = '  -0.34%  '
Set-TextValue 'D19' '349.63'
$ws.Range('E19').Value = '  +0.96%  '
Set-TextValue 'D20' '4.51'
$ws.Range('E20').Value = '  -2.40%  '
$ws.Range('E21').Value = '  -1.91%  '
$ws.Range('E22').Value = '  +0.69%  '
$ws.Range('E23').Value = '  +0.20%  '
Set-TextValue 'D24' '62.34'
$ws.Range('E24').Value = '  +2.35%  '
Set-TextValue 'D25' '0.416'
$ws.Range('E25').Value = '  -1.96%  '
$ws.Range('E26').Value = '  +2.81%  '
$ws.Range('E27').Value = '  +0.79%  '
Set-TextValue 'D28' '0.0₃0805'
$ws.Range('E28').Value = '  -2.30%  '
Set-TextValue 'D29' '7.16'
$ws.Range('E29').Value = '  -1.05%  '
$ws.Range('E30').Value = '  +0.10%  '
Set-TextValue 'D31' '6.50'
$ws.Range('E31').Value = '  -0.38%  '
$ws.Range('B32').Value = 'PancakeSwap'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue 'D32' '1.58'
$ws.Range('E32').Value = '  -0.20%  '
$ws.Range('B33').Value = 'EthereumClassic'
$ws.Range('C33').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue 'D33' '18.92'
$ws.Range('E33').Value = '  -0.59%  '
Set-TextValue 'D34' '149.65'
$ws.Range('E34').Value = '  +0.05%  '
Set-TextValue 'D35' '4.07'
$ws.Range('E35').Value = '  +0.65%  '
Set-TextValue 'D36' '0.944'
$ws.Range('E36').Value = '  -11.84%  '
$ws.Range('E37').Value = '  +0.87%  '
Set-TextValue 'D38' '0.864'
$ws.Range('E38').Value = '  -1.07%  '
Set-TextValue 'D39' '36.64'
$ws.Range('E39').Value = '  -0.35%  '
$ws.Range('E40').Value = '  +3.60%  '
Set-TextValue 'D41' '3.68'
$ws.Range('E41').Value = '  -1.98%  '
$ws.Range('B42').Value = 'Bittensor'
$ws.Range('C42').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue 'D42' '278.67'
$ws.Range('E42').Value = '  -2.05%  '
$ws.Range('B43').Value = 'Stellar'
$ws.Range('C43').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue 'D43' '0.0991'
$ws.Range('E43').Value = '  -0.54%  '
$ws.Range('E44').Value = '  +0.42%  '
Set-TextValue 'D45' '0.601'
$ws.Range('E45').Value = '  -2.95%  '
Set-TextValue 'D46' '19.58'
$ws.Range('E46').Value = '  -1.27%  '
Set-TextValue 'D47' '2.073.99'
$ws.Range('E47').Value = '  +4.13%  '
$ws.Range('E48').Value = '  -2.87%  '
Set-TextValue 'D49' '0.0232'
$ws.Range('E49').Value = '  -0.65%  '
Set-TextValue 'D50' '10.32'
$ws.Range('E50').Value = '  +0.47%  '
Set-TextValue 'D51' '4.72'
$ws.Range('E51').Value = '  -1.06%  '
